$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A47").Value = "Elia Barozzi"
$ws.Range("B47").Value = "Elia Barozzi | I Magnifici"
$ws.Range("C47").Value = "Carlo Stedile | MAI UNA GIOIA"
$ws.Range("D47").Value = "Luca Frasca | Clitoriders"
$ws.Range("E47").Value = "FEDERICO NICOLODI | U.S. Guarna"
$ws.Range("F47").Value = "Davide Fontana | SBARX"
